$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old sample data (win/loss/role/damage rows) but keep the
# existing cell formatting (border/bold/alignment) on the header + first
# data column cells.
$ws.Range("A1:C3").ClearContents()

# Extend the formatted first column down through row 5 for new game rows
# (started working on the "add game" button/flow).
$ws.Range("A2").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection where the author left it.
$ws.Range("H10").Select() | Out-Null
